$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 17 new rows (53-69) with the same column-A numbering style as existing rows
$ws.Range("A52").Copy()
$ws.Range("A53:A69").PasteSpecial(-4122)
for ($r = 53; $r -le 69; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update cell values per diff
$ws.Range("B2").Value = 'NSE:FDC'
$ws.Range("C2").Value = 'NSE:ADFFOODS'
$ws.Range("D2").Value = 'NSE:MAXHEALTH'
$ws.Range("E2").Value = 'NSE:BAJFINANCE'
$ws.Range("F2").Value = 'NSE:MUTHOOTFIN'
$ws.Range("B3").Value = 'NSE:GRMOVER'
$ws.Range("C3").Value = 'NSE:AGI'
$ws.Range("E3").Value = 'NSE:CUMMINSIND'
$ws.Range("B4").Value = 'NSE:LICMFGOLD'
$ws.Range("C4").Value = 'NSE:AGRITECH'
$ws.Range("E4").ClearContents()
$ws.Range("B5").Value = 'NSE:MALLCOM'
$ws.Range("C5").Value = 'NSE:ALANKIT'
$ws.Range("E5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = 'NSE:ALOKINDS'
$ws.Range("E6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 'NSE:APOLLOTYRE'
$ws.Range("E7").ClearContents()
$ws.Range("C8").Value = 'NSE:ARIES'
$ws.Range("E8").ClearContents()
$ws.Range("C9").Value = 'NSE:ARTEMISMED'
$ws.Range("E9").ClearContents()
$ws.Range("C10").Value = 'NSE:BAJAJCON'
$ws.Range("E10").ClearContents()
$ws.Range("C11").Value = 'NSE:BAJAJHIND'
$ws.Range("E11").ClearContents()
$ws.Range("C12").Value = 'NSE:BANARISUG'
$ws.Range("E12").ClearContents()
$ws.Range("C13").Value = 'NSE:BASML'
$ws.Range("E13").ClearContents()
$ws.Range("C14").Value = 'NSE:BBTC'
$ws.Range("C15").Value = 'NSE:BCLIND'
$ws.Range("C16").Value = 'NSE:BFUTILITIE'
$ws.Range("C17").Value = 'NSE:BIOFILCHEM'
$ws.Range("C18").Value = 'NSE:BIRLAMONEY'
$ws.Range("C19").Value = 'NSE:BSE'
$ws.Range("C20").Value = 'NSE:CORDSCABLE'
$ws.Range("C21").Value = 'NSE:DCW'
$ws.Range("C22").Value = 'NSE:DELTACORP'
$ws.Range("C23").Value = 'NSE:DEVIT'
$ws.Range("C24").Value = 'NSE:DNAMEDIA'
$ws.Range("C25").Value = 'NSE:DOLLAR'
$ws.Range("C26").Value = 'NSE:DWARKESH'
$ws.Range("C27").Value = 'NSE:ELDEHSG'
$ws.Range("C28").Value = 'NSE:FCL'
$ws.Range("C29").Value = 'NSE:FCSSOFT'
$ws.Range("C30").Value = 'NSE:GAEL'
$ws.Range("C31").Value = 'NSE:GHCL'
$ws.Range("C32").Value = 'NSE:GIPCL'
$ws.Range("C33").Value = 'NSE:GREENPOWER'
$ws.Range("C34").Value = 'NSE:GRINFRA'
$ws.Range("C35").Value = 'NSE:GSS'
$ws.Range("C36").Value = 'NSE:HARRMALAYA'
$ws.Range("C37").Value = 'NSE:HARSHA'
$ws.Range("C38").Value = 'NSE:HITECH'
$ws.Range("C39").Value = 'NSE:HMAAGRO'
$ws.Range("C40").Value = 'NSE:INDOWIND'
$ws.Range("C41").Value = 'NSE:IRB'
$ws.Range("C42").Value = 'NSE:JAGSNPHARM'
$ws.Range("C43").Value = 'NSE:JINDALSAW'
$ws.Range("C44").Value = 'NSE:JKIL'
$ws.Range("C45").Value = 'NSE:JMFINANCIL'
$ws.Range("C46").Value = 'NSE:KMSUGAR'
$ws.Range("C47").Value = 'NSE:KOTARISUG'
$ws.Range("C48").Value = 'NSE:LANDMARK'
$ws.Range("C49").Value = 'NSE:MAHSEAMLES'
$ws.Range("C50").Value = 'NSE:MAPMYINDIA'
$ws.Range("C51").Value = 'NSE:MMFL'
$ws.Range("C52").Value = 'NSE:MOIL'
$ws.Range("C53").Value = 'NSE:MUFIN'
$ws.Range("C54").Value = 'NSE:NATCOPHARM'
$ws.Range("C55").Value = 'NSE:NIITLTD'
$ws.Range("C56").Value = 'NSE:NSLNISP'
$ws.Range("C57").Value = 'NSE:ONWARDTEC'
$ws.Range("C58").Value = 'NSE:OSIAHYPER'
$ws.Range("C59").Value = 'NSE:PANACEABIO'
$ws.Range("C60").Value = 'NSE:PARADEEP'
$ws.Range("C61").Value = 'NSE:PATELENG'
$ws.Range("C62").Value = 'NSE:PFOCUS'
$ws.Range("C63").Value = 'NSE:PFS'
$ws.Range("C64").Value = 'NSE:PRECWIRE'
$ws.Range("C65").Value = 'NSE:PRICOLLTD'
$ws.Range("C66").Value = 'NSE:REDTAPE'
$ws.Range("C67").Value = 'NSE:REFEX'
$ws.Range("C68").Value = 'NSE:RICOAUTO'
$ws.Range("C69").Value = 'NSE:RSWM'

$ws.Range("A1").Select()
